$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 96

$ws.Cells.Item($newRow, 1).Value = 3
$ws.Cells.Item($newRow, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"

$ws.Cells.Item(95, 4).Copy()
$ws.Cells.Item($newRow, 4).PasteSpecial(-4122)
$ws.Cells.Item($newRow, 4).Value = (Get-Date -Year 2022 -Month 2 -Day 3 -Hour 0 -Minute 0 -Second 0)

$ws.Cells.Item($newRow, 5).Value = 5
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100107
$ws.Cells.Item($newRow, 8).Value = "Otros"
$ws.Cells.Item($newRow, 9).Value = 100107011
$ws.Cells.Item($newRow, 10).Value = "Tuna"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 50
$ws.Cells.Item($newRow, 14).Value = 20000
$ws.Cells.Item($newRow, 15).Value = 20000
$ws.Cells.Item($newRow, 16).Value = 20000
$ws.Cells.Item($newRow, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item($newRow, 18).Value = "Provincia de Limarí"
$ws.Cells.Item($newRow, 19).Value = 1000
$ws.Cells.Item($newRow, 20).Value = 20
